# SqlServer.xlsx -- "unify the conception of DataNode, DataTable, Entity."
#
# The worksheet that used to be a generic "Property1" sheet is renamed to
# "DataNode" as part of a repo-wide renaming pass, and the workbook's saved
# selection moves from A9 to E23 (the cell the author was last looking at
# when the file was saved).

$wb = $excel.ActiveWorkbook

# Rename the (only) worksheet from "Property1" to "DataNode".
$ws = $wb.Worksheets.Item("Property1")
$ws.Name = "DataNode"

# Make sure it's the active sheet, then move/save the selection to E23
# (previously A9), matching the cursor position recorded in the file.
$ws.Activate()
$ws.Range("E23").Select()
